$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "PizzaCluster is…….." -> real project-summary sentence.
# The match/replacement text is identical on both sides of the ellipsis
# run boundary, so Find/Execute rewrites the single run in place.
# ------------------------------------------------------------------
$find1 = $d.Content.Find
[void]$find1.Execute(" is……..", $false, $false, $false, $false, $false, $true, 1, $false, " is a server side database retrieving system. It allows users to browse pizza restaurant and ordering pizza remotely using any web browsing device. It tracks user information and activity histories to intense the experiences with online pizza ordering. ", 2)

# ------------------------------------------------------------------
# Change 2: merge the lone " " run with the following "Once we
# assigned..." run into a single run (leave the preceding "...epic."
# run untouched). Replacing "Once we assigned..." with itself forces
# the engine to re-flow just the runs the match touches.
# ------------------------------------------------------------------
$find2 = $d.Content.Find
[void]$find2.Execute("Once we assigned our tasks, we started to meet after every class to talk about what we had done since the previous meeting and to see if anyone had anything important to share. This was our version of scrum meetings.", $false, $false, $false, $false, $false, $true, 1, $false, "Once we assigned our tasks, we started to meet after every class to talk about what we had done since the previous meeting and to see if anyone had anything important to share. This was our version of scrum meetings.", 2)

# ------------------------------------------------------------------
# Change 3: merge "...instead of" + " dividing tasks..." into one run.
# ------------------------------------------------------------------
$find3 = $d.Content.Find
[void]$find3.Execute(" dividing tasks which again crippled our progress more than significantly.", $false, $false, $false, $false, $false, $true, 1, $false, " dividing tasks which again crippled our progress more than significantly.", 2)

# ------------------------------------------------------------------
# Change 4: "Iteration-1 " + "What" + " is done, what is to be done?"
# (3 runs, 2 proofErr tags) -> single run, no proofErr markers.
# ------------------------------------------------------------------
$find4 = $d.Content.Find
[void]$find4.Execute("Iteration-1 What is done, what is to be done?", $false, $false, $false, $false, $false, $true, 1, $false, "Iteration-1 What is done, what is to be done?", 2)

# ------------------------------------------------------------------
# Change 5: turn the blank paragraph that used to follow the
# "Iteration-1 What is done..." heading into four new bullet-style
# paragraphs (720-twip first-line indent), then relocate the hidden
# "_GoBack" bookmark from the old tab paragraph onto the end of the
# final "Basic GUI." paragraph.
# ------------------------------------------------------------------
$headerPara = $d.Paragraphs.Item(15)
$newItems = @("Basic structure of the website.", "Basic user database system.", "Demo pizza restaurant database.", "Basic GUI.")

$insertAfter = $headerPara
foreach ($itemText in $newItems) {
    $insertAfter.Range.InsertParagraphAfter()
    $insertAfter = $insertAfter.Next()
    $insertAfter.Range.Text = $itemText
    $insertAfter.Format.FirstLineIndent = 36
}

# $insertAfter is now the "Basic GUI." paragraph. Append a throwaway
# placeholder character so the bookmark's anchor position isn't the
# exact "end of paragraph content" boundary, add/relocate the
# "_GoBack" bookmark there, then delete the placeholder again - the
# bookmark stays put because Word bookmarks track surrounding edits.
$insertAfter.Range.InsertAfter("X")
$guiPara = $d.Paragraphs.Item($insertAfter.Index)
$placeholderPos = $guiPara.Range.End - 2
$bookmarkRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()

Write-Output "edit complete"
